$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 7
$ws.Range("D2").Value = 44567
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 2400
$ws.Range("O2").Value = 2400
$ws.Range("P2").Value = 2400
$ws.Range("R2").Value = "Región de La Araucanía"
$ws.Range("S2").Value = 2400

# Row 3 <- old row 2
$ws.Range("D3").Value = 44214
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 1800
$ws.Range("O3").Value = 1800
$ws.Range("P3").Value = 1800
$ws.Range("R3").Value = "Región de La Araucanía"
$ws.Range("S3").Value = 1800

# Row 4 <- old row 10
$ws.Range("D4").Value = 44616
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 3200
$ws.Range("O4").Value = 3200
$ws.Range("P4").Value = 3200
$ws.Range("R4").Value = "Región de La Araucanía"
$ws.Range("S4").Value = 3200

# Row 5 <- old row 6
$ws.Range("D5").Value = 44175
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 5000

# Row 6 <- old row 3
$ws.Range("D6").Value = 44176
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 3000
$ws.Range("O6").Value = 3000
$ws.Range("P6").Value = 3000
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 3000

# Row 7 <- old row 8
$ws.Range("D7").Value = 44551
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 4500
$ws.Range("O7").Value = 4500
$ws.Range("P7").Value = 4500
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 4500

# Row 8 <- old row 11
$ws.Range("D8").Value = 44215
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 2800
$ws.Range("O8").Value = 2800
$ws.Range("P8").Value = 2800
$ws.Range("R8").Value = "Región de La Araucanía"
$ws.Range("S8").Value = 2800

# Row 9 <- old row 4
$ws.Range("D9").Value = 44592
$ws.Range("M9").Value = 5
$ws.Range("N9").Value = 7500
$ws.Range("O9").Value = 7500
$ws.Range("P9").Value = 7500
$ws.Range("R9").Value = "Región de La Araucanía"
$ws.Range("S9").Value = 7500

# Row 10 <- old row 5
$ws.Range("D10").Value = 44323
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 3200
$ws.Range("O10").Value = 3200
$ws.Range("P10").Value = 3200
$ws.Range("R10").Value = "Región de La Araucanía"
$ws.Range("S10").Value = 3200

# Row 11 <- old row 9
$ws.Range("D11").Value = 44574
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 3000
$ws.Range("O11").Value = 3000
$ws.Range("P11").Value = 3000
$ws.Range("R11").Value = "Región de La Araucanía"
$ws.Range("S11").Value = 3000
